# Apply the commit "add abstract layers of components, refactoring, add glow,
# polish particles, change resolution" edits to the dev-notes workbook.
#
# Changes:
#  1. TODO Before 0.0.1: update two TODO item descriptions to mention STAMINA.
#  2. TODO Before 0.0.1: move the active selection (no longer the active tab).
#  3. Logs: append three new dated log entries (rows 35-37).
#  4. Logs: becomes the active sheet/tab, with a new selection.

$wb = $excel.ActiveWorkbook

# --- 1. Update TODO item text on "TODO Before 0.0.1" ---------------------
$wsTodo = $wb.Worksheets.Item("TODO Before 0.0.1")
$wsTodo.Range("B12").Value = "create enemy mp, hp, xp system AND STAMINA"
$wsTodo.Range("B7").Value  = "enrich enemies and character with basic HP, MP AND STAMINA"

# --- 3. Append new log rows to "Logs" -------------------------------------
$wsLogs = $wb.Worksheets.Item("Logs")

# Row 35 - copy formatting from the row above (row 34) then set new values.
$wsLogs.Range("A34:B34").Copy()
$wsLogs.Range("A35:B35").PasteSpecial(-4122)
$wsLogs.Range("A35").Value = 45460
$wsLogs.Range("B35").Value = "Add stamina component, fix inheritence for health and mana components"

# Row 36 - same single-line formatting as row 35/34.
$wsLogs.Range("A34:B34").Copy()
$wsLogs.Range("A36:B36").PasteSpecial(-4122)
$wsLogs.Range("A36").Value = 45462
$wsLogs.Range("B36").Value = "move player skills to controllers (including jump, run, aim etc.) "

# Row 37 - copy formatting from a previously-wrapped 2-line row (row 31).
$wsLogs.Range("A31:B31").Copy()
$wsLogs.Range("A37:B37").PasteSpecial(-4122)
$wsLogs.Range("A37").Value = 45463
$wsLogs.Range("B37").Value = "refactor again - basic stamina and jump controllers used from userSkillController now, fix bugs when skill consume stamina continuasly, add lock system for mana and stamina skills - to avoid paralel skill usages"
$wsLogs.Rows.Item(37).RowHeight = 28.8

$excel.CutCopyMode = $false

# --- 2/4. Update sheet selections and switch the active tab to Logs ------
$wsTodo.Activate()
$wsTodo.Range("B14").Select()

$wsLogs.Activate()
$wsLogs.Range("B39").Select()
